$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Extend the table by three rows (65-67), reusing the formatting
# (borders/bold for column A, datetime number format for column E) from
# the last existing data row (64) so the new rows are styled identically
# to the rest of the sheet.
$ws.Range("A64:V64").Copy() | Out-Null
$ws.Range("A65:V67").PasteSpecial(-4122) | Out-Null

# Row 65: FUS Rabat 1-1 Berkane
$ws.Range("A65").Value = 64
$ws.Range("B65").Value = "morocco"
$ws.Range("C65").Value = "botola-pro"
$ws.Range("D65").Value = "2023-2024"
$ws.Range("E65").Value = 45241.66666666666
$ws.Range("F65").Value = "FUS Rabat"
$ws.Range("G65").Value = 1
$ws.Range("H65").Value = "Berkane"
$ws.Range("I65").Value = 1
$ws.Range("J65").Value = 2.31
$ws.Range("K65").Value = "10/11/2023 04:13"
$ws.Range("L65").Value = 2.66
$ws.Range("M65").Value = "11/11/2023 15:58"
$ws.Range("N65").Value = 2.75
$ws.Range("O65").Value = "10/11/2023 04:13"
$ws.Range("P65").Value = 2.52
$ws.Range("Q65").Value = "11/11/2023 15:58"
$ws.Range("R65").Value = 3.16
$ws.Range("S65").Value = "10/11/2023 04:13"
$ws.Range("T65").Value = 3.35
$ws.Range("U65").Value = "11/11/2023 15:44"
$ws.Range("V65").Value = "https://www.betexplorer.com/football/morocco/botola-pro/fus-rabat-berkane/xjqoUUxM/"

# Row 66: Jeunesse Sportive Soualem 0-0 Raja Casablanca
$ws.Range("A66").Value = 65
$ws.Range("B66").Value = "morocco"
$ws.Range("C66").Value = "botola-pro"
$ws.Range("D66").Value = "2023-2024"
$ws.Range("E66").Value = 45241.76041666666
$ws.Range("F66").Value = "Jeunesse Sportive Soualem"
$ws.Range("G66").Value = 0
$ws.Range("H66").Value = "Raja Casablanca"
$ws.Range("I66").Value = 0
$ws.Range("J66").Value = 4.8
$ws.Range("K66").Value = "10/11/2023 06:42"
$ws.Range("L66").Value = 6.83
$ws.Range("M66").Value = "11/11/2023 18:11"
$ws.Range("N66").Value = 3.39
$ws.Range("O66").Value = "10/11/2023 06:42"
$ws.Range("P66").Value = 3.64
$ws.Range("Q66").Value = "11/11/2023 18:11"
$ws.Range("R66").Value = 1.63
$ws.Range("S66").Value = "10/11/2023 06:42"
$ws.Range("T66").Value = 1.54
$ws.Range("U66").Value = "11/11/2023 18:11"
$ws.Range("V66").Value = "https://www.betexplorer.com/football/morocco/botola-pro/jeunesse-sportive-soualem-raja-casablanca/GtrkTlhS/"

# Row 67: Olympique de Safi 1-0 FAR Rabat
$ws.Range("A67").Value = 66
$ws.Range("B67").Value = "morocco"
$ws.Range("C67").Value = "botola-pro"
$ws.Range("D67").Value = "2023-2024"
$ws.Range("E67").Value = 45241.85416666666
$ws.Range("F67").Value = "Olympique de Safi"
$ws.Range("G67").Value = 1
$ws.Range("H67").Value = "FAR Rabat"
$ws.Range("I67").Value = 0
$ws.Range("J67").Value = 3.92
$ws.Range("K67").Value = "10/11/2023 08:42"
$ws.Range("L67").Value = 3.52
$ws.Range("M67").Value = "11/11/2023 20:18"
$ws.Range("N67").Value = 3
$ws.Range("O67").Value = "10/11/2023 08:42"
$ws.Range("P67").Value = 2.75
$ws.Range("Q67").Value = "11/11/2023 20:18"
$ws.Range("R67").Value = 1.93
$ws.Range("S67").Value = "10/11/2023 08:42"
$ws.Range("T67").Value = 2.36
$ws.Range("U67").Value = "11/11/2023 20:18"
$ws.Range("V67").Value = "https://www.betexplorer.com/football/morocco/botola-pro/olympique-de-safi-far-rabat/nPUGO80q/"
